$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.345.14"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.526.39"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.37"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.59"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "2.527.03"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "2.974.46"
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").Value = "58.494.52"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.16"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "2.524.70"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.69"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "322.19"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +8.10%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.74"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.161"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "167.94"
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.20"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.15"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.43"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.46"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.35"
$ws.Range("E44").Value = "  +5.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.00"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0503"
$ws.Range("E48").Value = "  +3.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.81"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.93"
$ws.Range("E51").Value = "  +0.65%  "
